$d = $word.ActiveDocument

$d.Content.Find.Execute("Rapport de Solde de Tout Compte pour le Logement ID: 1", $true, $false, $false, $false, $false, $true, 1, $false, "Rapport de Solde de Tout Compte pour le Logement ID: 121", 2)
$d.Content.Find.Execute("Prix des consommations : 253.4", $true, $false, $false, $false, $false, $true, 1, $false, "Prix des consommations : 0.0", 2)
$d.Content.Find.Execute("Loyers impayés : 27000.0", $true, $false, $false, $false, $false, $true, 1, $false, "Loyers impayés : 20250.0", 2)
$d.Content.Find.Execute("Charges récupérables : 50.0", $true, $false, $false, $false, $false, $true, 1, $false, "Charges récupérables : 0.0", 2)
